# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 11621
    5  = 1057
    6  = 125
    10 = 10990
    11 = 4230
    12 = 22
    18 = 9
    19 = 147
    21 = 11180
    22 = 10999
    24 = 34
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
